# Auto-generated edit script: updates cryptocurrency Price (D) and Volume(1h) (E)
# columns on Sheet1 to match the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.222.32"
$ws.Range("E2").Value = "  -1.63%  "
$ws.Range("D3").Value = "2.513.31"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'572.19"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "'166.71"
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +1.65%  "
$ws.Range("D9").Value = "2.509.20"
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("E10").Value = "  +0.30%  "
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").Value = "'0.353"
$ws.Range("E12").Value = "  +3.38%  "
$ws.Range("D13").Value = "'4.92"
$ws.Range("E13").Value = "  +2.80%  "
$ws.Range("D14").Value = "2.975.27"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "69.063.19"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("D16").Value = "'0.0000175"
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("D17").Value = "'24.80"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "2.514.12"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "'11.35"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").Value = "'7.60"
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("D21").Value = "'348.68"
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").Value = "  +1.00%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "'70.25"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("D26").Value = "'3.98"
$ws.Range("E26").Value = "  -1.40%  "
$ws.Range("D27").Value = "'8.95"
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("D28").Value = "2.644.62"
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -4.08%  "
$ws.Range("D30").Value = "0.0₃0892"
$ws.Range("E30").Value = "  -1.33%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").Value = "'462.45"
$ws.Range("E32").Value = "  -2.87%  "
$ws.Range("E33").Value = "  -3.93%  "
$ws.Range("E34").Value = "  -1.77%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("D37").Value = "'156.48"
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("D38").Value = "'19.00"
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("D39").Value = "'18.53"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "'4.74"
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("D44").Value = "'38.02"
$ws.Range("E45").Value = "  -13.28%  "
$ws.Range("D46").Value = "'2.26"
$ws.Range("E46").Value = "  -4.88%  "
$ws.Range("D47").Value = "'141.67"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("D50").Value = "'0.0730"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("E51").Value = "  -2.80%  "
